$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the duplicated "aaa" name in A3 to "aaa1" to make it distinct
# (subset used to remove duplicates filtered by email)
$ws.Range("A3").Value = "aaa1"

# Move the selection to A3 to match the resulting file's cursor position
$ws.Range("A3").Select()
